# "tidying up spm excel sheets"
#
# - Add a new "info" sheet at the front of the workbook with some metadata
#   about the dataset (author, last update, code link, units).
# - Rename the "fraction" header used on each of the existing data sheets
#   to "shares" (values are unchanged).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "info" sheet as the very first tab.
# ------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$info = $wb.Worksheets.Add($firstSheet)
$info.Name = "info"

$info.Range("A1").Value = "Author"
$info.Range("B1").Value = "William F. Lamb"

$info.Range("A2").Value = "Last update"
$info.Range("B2").Value = "2021-10-14 10:09:23"

$info.Range("A3").Value = "Code"
$info.Range("B3").Value = "https://github.com/mcc-apsis/AR6-Emissions-trends-and-drivers/blob/master/R/Analysis%20and%20figures/direct_indirect_emissions.Rmd"

$info.Range("A4").Value = ""
$info.Range("B4").Value = ""

$info.Range("A5").Value = "Units"
$info.Range("B5").Value = ""

$info.Range("A6").Value = "GHG"
$info.Range("B6").Value = "GtCO2eq"

$info.Range("A7").Value = "shares"
$info.Range("B7").Value = "%"

# ------------------------------------------------------------------
# 2) Rename the "fraction" column header to "shares" on each of the
#    pre-existing data sheets (the data itself does not change).
# ------------------------------------------------------------------
$direct = $wb.Worksheets.Item("direct emissions")
$direct.Range("C1").Value = "shares"

$indirect = $wb.Worksheets.Item("indirect emissions")
$indirect.Range("E1").Value = "shares"

$subsectors = $wb.Worksheets.Item("indirect emissions - subsectors")
$subsectors.Range("E1").Value = "shares"
